# Fix a transaction that was mis-attributed to "Sales" instead of
# "Production" after the item's name changed mid-edit. Moves the
# -66150 amount from Sales (Week 2) to Production (Week 2) in both the
# "Total Cost per Section" summary table and the weekly breakdown
# table, and recomputes the affected row/column totals.
#
# NOTE: table/cell handles are re-fetched via $d.Tables.Item(...) for
# every single write (rather than cached in a variable and reused)
# because writing through a stale cached Table/Cell handle after the
# document has already been mutated once can desynchronize it from the
# live document and corrupt unrelated rows.

$d = $word.ActiveDocument

# --- Table 2: "I. Total Cost per Section" summary -------------------
# b. Production -> Total column: 0 -> -66150
$d.Tables.Item(2).Cell(6, 3).Range.Text = "-66150"

# g. Sales -> Total column: -105 -> 0
$d.Tables.Item(2).Cell(11, 3).Range.Text = "0"

# Grand Total row -> Total column: -105 -> -66150
$d.Tables.Item(2).Cell(13, 3).Range.Text = "-66150"

# --- Table 3: weekly breakdown ---------------------------------------
# Production row: WEEK 2 and Total columns: 0 -> -66150
$d.Tables.Item(3).Cell(3, 3).Range.Text = "-66150"
$d.Tables.Item(3).Cell(3, 7).Range.Text = "-66150"

# Sales row: WEEK 2 and Total columns: -105 -> 0
$d.Tables.Item(3).Cell(8, 3).Range.Text = "0"
$d.Tables.Item(3).Cell(8, 7).Range.Text = "0"

# TOTAL row: WEEK 2 and Total columns: -105 -> -66150
$d.Tables.Item(3).Cell(10, 3).Range.Text = "-66150"
$d.Tables.Item(3).Cell(10, 7).Range.Text = "-66150"
